$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert two new rows above the existing "OrchestratorQueueFolder" row (row 3),
# pushing everything below it down by two rows. Excel carries the formatting
# of the row above (row 2) down into the newly inserted rows.
$ws.Rows("3:4").Insert()
$ws.Rows("3:4").RowHeight = 14.25

# Fill in the two newly inserted rows with the new config entries, and update
# the queue-name value (row 2) to the new process name.
$ws.Range("B3").Value = "Data\Input\CandidatesInfo.xlsx"
$ws.Range("B2").Value = "GenerateOfferLetter"
$ws.Range("A3").Value = "CadidatesInfoFilePath"

$ws.Range("A4").Value = "CadidatesInfoSheetName"
$ws.Range("B4").Value = "OutputData"

# Keep the selection where the new content was entered, matching the saved file.
[void]$ws.Range("B4").Select()
